$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New abbreviation / complete-word rows appended to the glossary table.
$newRows = @(
    @("B/R",   "back reaming"),
    @("std",   "tool joint"),
    @("ID ",   "Inner diameter"),
    @("OD",    "Out diamter"),
    @("SLB",   "Schlumberger"),
    @("HAL",   "Halliburton"),
    @("BHGE",  "Baker Hughes"),
    @("WFD",   "Weatherford"),
    @("WL",    "Wireline logging"),
    @("STDS",  "Stands"),
    @("TJ",    "tool joint"),
    @("CO",    "Cross over"),
    @("XMAS",  "Christmas tree"),
    @("LCM",   "Lost circulation material"),
    @("PDS",   "Product data sheet"),
    @("SDS",   "Safety data sheet"),
    @("MSDS",  "Material safety data sheet"),
    @("btm",   "bottom")
)

$startRow = 42
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

$lastRow = $startRow + $newRows.Count - 1

# Grow the "Table1" structured table to cover the newly-added rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B$lastRow"))

# Column B widened (longer glossary entries like "Material safety data sheet").
$ws.Columns.Item(2).ColumnWidth = 23

# Leave the selection where the author left it after typing the new rows.
$ws.Range("B64").Select()
